# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail table (B15:J42) lists, for every mora
# period, each worker's doc type / doc number / name / period / value.
# This edit re-sorts the table so it is grouped by period (oldest -> newest:
# 2008, 2009, 2010, 2011, 2012, 2101, 2102, 2103, 2104) and, within each
# period, lists the three workers in the same fixed order. The "oldest"
# value (32586) that used to be attached to the first-listed period (2104)
# now travels with the period that ends up last (2104 again, since it is
# still the most recent period) - i.e. rows 40-42 keep the 32586 value
# while every other row uses 37600.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2008", "2009", "2010", "2011", "2012", "2101", "2102", "2103", "2104")
$workers = @(
    @("CC", "73557907", "ROBERT GARCES COTA"),
    @("CC", "73132230", "JOSE FREDYS GELIZ PEREZ"),
    @("CC", "11051258", "NELSON MIGUEL CALDERA RICARDO")
)

$row = 16
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        $tipoDoc = $worker[0]
        $nroDoc = $worker[1]
        $nombre = $worker[2]

        $valorMora = 37600
        if ($period -eq "2104") {
            $valorMora = 32586
        }

        $ws.Cells.Item($row, 2).Value = $tipoDoc
        $ws.Cells.Item($row, 3).Value = $nroDoc
        $ws.Cells.Item($row, 4).Value = $nombre
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = $valorMora

        $row = $row + 1
    }
}
